$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 3377.3157
$ws.Range("I98").Value = 3869.2144
$ws.Range("J98").Value = 2000
$ws.Range("K98").Value = 3869.2144
$ws.Range("L98").Value = 2000
$ws.Range("M98").Value = -2371.2144
$ws.Range("N98").Value = -4996
$ws.Range("H122").Value = 3377.3157
$ws.Range("I122").Value = 3869.2144
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 11607.6432
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -9157.643199999999
$ws.Range("N122").Value = -10900
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").Value = ""
$ws.Range("H133").Value = 55000
$ws.Range("J133").Value = 55000
$ws.Range("L133").Value = 55000
$ws.Range("N133").Value = -65120
$ws.Range("H137").Value = 428147.88
$ws.Range("I137").Value = 679509.75
$ws.Range("J137").Value = 2766.3076
$ws.Range("K137").Value = 2038529.25
$ws.Range("L137").Value = 8298.9228
$ws.Range("M137").Value = -2035979.25
$ws.Range("N137").Value = -13398.9228

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3112.45
$ws.Range("I61").Value = 2394.3572
$ws.Range("J61").Value = 4788
$ws.Range("K61").Value = 2394.3572
$ws.Range("L61").Value = 4788
$ws.Range("M61").Value = -2182.3572
$ws.Range("N61").Value = -5212
$ws.Range("H62").Value = 49800
$ws.Range("J62").Value = 49800
$ws.Range("L62").Value = 49800
$ws.Range("N62").Value = -51048
$ws.Range("H65").Value = 49800
$ws.Range("J65").Value = 49800
$ws.Range("L65").Value = 149400
$ws.Range("N65").Value = -155640
$ws.Range("H122").Value = 1541.45
$ws.Range("I122").Value = 1594.5385
$ws.Range("J122").Value = 1442.8572
$ws.Range("K122").Value = 4783.6155
$ws.Range("L122").Value = 4328.571599999999
$ws.Range("M122").Value = -2333.6155
$ws.Range("N122").Value = -9228.571599999999
$ws.Range("H132").Value = 1439615.1
$ws.Range("I132").Value = 1918209.1
$ws.Range("K132").Value = 5754627.300000001
$ws.Range("M132").Value = -5752097.300000001
$ws.Range("H136").Value = 3112.45
$ws.Range("I136").Value = 2394.3572
$ws.Range("J136").Value = 4788
$ws.Range("K136").Value = 7183.071599999999
$ws.Range("L136").Value = 14364
$ws.Range("M136").Value = -4633.071599999999
$ws.Range("N136").Value = -19464

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4833.3335
$ws.Range("I105").Value = 6000
$ws.Range("J105").Value = 2500
$ws.Range("K105").Value = 6000
$ws.Range("L105").Value = 2500
$ws.Range("M105").Value = -4253
$ws.Range("N105").Value = -5994
$ws.Range("H107").Value = 224792.25
$ws.Range("I107").Value = 307315.1
$ws.Range("J107").Value = 1495.1177
$ws.Range("K107").Value = 307315.1
$ws.Range("L107").Value = 1495.1177
$ws.Range("M107").Value = -305395.1
$ws.Range("N107").Value = -5335.1177
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").Value = ""
$ws.Range("H135").Value = 79900
$ws.Range("J135").Value = 79900
$ws.Range("L135").Value = 79900
$ws.Range("N135").Value = -90040

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1190.9117
$ws.Range("I31").Value = 980.0357
$ws.Range("J31").Value = 2175
$ws.Range("K31").Value = 980.0357
$ws.Range("L31").Value = 2175
$ws.Range("M31").Value = -685.0357
$ws.Range("N31").Value = -2765
$ws.Range("H34").Value = 1190.9117
$ws.Range("I34").Value = 980.0357
$ws.Range("J34").Value = 2175
$ws.Range("K34").Value = 980.0357
$ws.Range("L34").Value = 2175
$ws.Range("M34").Value = -778.0357
$ws.Range("N34").Value = -2579
$ws.Range("H35").Value = 972.7692
$ws.Range("I35").Value = 972.7692
$ws.Range("K35").Value = 972.7692
$ws.Range("M35").Value = -678.7692

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 580.46875
$ws.Range("J113").Value = 569.04
$ws.Range("L113").Value = 1707.12
$ws.Range("N113").Value = -6047.12
$ws.Range("H117").Value = 2479.2307
$ws.Range("I117").Value = 1060
$ws.Range("J117").Value = 3366.25
$ws.Range("K117").Value = 3180
$ws.Range("L117").Value = 10098.75
$ws.Range("M117").Value = 262
$ws.Range("N117").Value = -16982.75
$ws.Range("H129").Value = 587.8
$ws.Range("I129").Value = 547.25
$ws.Range("J129").Value = 750
$ws.Range("K129").Value = 1641.75
$ws.Range("L129").Value = 2250
$ws.Range("M129").Value = 3358.25
$ws.Range("N129").Value = -12250

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").Value = ""
$ws.Range("H122").Value = 4187.5
$ws.Range("I122").Value = 6375
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 19125
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -16675
$ws.Range("N122").Value = -10900

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1879.4517
$ws.Range("I7").Value = 1778.091
$ws.Range("J7").Value = 2127.2222
$ws.Range("K7").Value = 1778.091
$ws.Range("L7").Value = 2127.2222
$ws.Range("M7").Value = -1666.091
$ws.Range("N7").Value = -2351.2222
$ws.Range("H40").Value = 1431.8276
$ws.Range("I40").Value = 1271.7916
$ws.Range("J40").Value = 2200
$ws.Range("K40").Value = 1271.7916
$ws.Range("L40").Value = 2200
$ws.Range("M40").Value = -1135.7916
$ws.Range("N40").Value = -2472
$ws.Range("H61").Value = 1000
$ws.Range("I61").Value = 833.3333
$ws.Range("K61").Value = 833.3333
$ws.Range("M61").Value = -631.3333
$ws.Range("H62").Value = 24900
$ws.Range("I62").Value = 20000
$ws.Range("K62").Value = 20000
$ws.Range("M62").Value = -19376
$ws.Range("H64").Value = 30075
$ws.Range("J64").Value = 30075
$ws.Range("L64").Value = 30075
$ws.Range("N64").Value = -30525
$ws.Range("H65").Value = 24900
$ws.Range("I65").Value = 20000
$ws.Range("K65").Value = 60000
$ws.Range("M65").Value = -56880
$ws.Range("H67").Value = 30075
$ws.Range("J67").Value = 30075
$ws.Range("L67").Value = 30075
$ws.Range("N67").Value = -31635
$ws.Range("H113").Value = 1000
$ws.Range("I113").Value = 833.3333
$ws.Range("K113").Value = 833.3333
$ws.Range("M113").Value = 1336.6667
$ws.Range("H126").Value = 1879.4517
$ws.Range("I126").Value = 1778.091
$ws.Range("J126").Value = 2127.2222
$ws.Range("K126").Value = 5334.272999999999
$ws.Range("L126").Value = 6381.6666
$ws.Range("M126").Value = -2864.272999999999
$ws.Range("N126").Value = -11321.6666
